# Apply 1% winsorization results for ln_carbon_intensity dependent variable
# Updated TWFE regression coefficients, std errors, t-values and p-values.
# (PowerShell's parser here has no literal support for scientific notation,
#  e.g. "1.23E-15", so very small/large magnitudes are expressed as a
#  division/multiplication so the resulting double matches exactly.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (const)
$ws.Range("B2").Value = -1.027810138928968 / 1000000000000000
$ws.Range("C2").Value = 0.002482392527945528
$ws.Range("D2").Value = -414.0401356185205 / 1000000000000000
$ws.Range("E2").Value = 0.9999999999996698

# Row 3 (DID)
$ws.Range("B3").Value = 0.05532890472139857
$ws.Range("C3").Value = 0.01053515919879201
$ws.Range("D3").Value = 5.251833757551831
$ws.Range("E3").Value = 1.602573171943789 / 10000000

# Row 4 (ln_real_gdp)
$ws.Range("B4").Value = -0.8092585493269315
$ws.Range("C4").Value = 0.04222877417436777
$ws.Range("D4").Value = -19.16367607511893
$ws.Range("E4").Value = 0

# Row 5 (ln_人口密度)
$ws.Range("B5").Value = -0.0188714528127612
$ws.Range("C5").Value = 0.03672625470189612
$ws.Range("D5").Value = -0.5138409284022881
$ws.Range("E5").Value = 0.607398056844755

# Row 6 (ln_金融发展水平)
$ws.Range("B6").Value = 0.01107475954305597
$ws.Range("C6").Value = 0.01992556504053506
$ws.Range("D6").Value = 0.5558065490502437
$ws.Range("E6").Value = 0.5783812727909732

# Row 7 (第二产业占GDP比重)
$ws.Range("B7").Value = 0.2940073584038606
$ws.Range("C7").Value = 0.06614079780070675
$ws.Range("D7").Value = 4.445174055652515
$ws.Range("E7").Value = 9.074980176437819 / 1000000
